# Updated cryptos list on Wed Sep  6 02:55:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.907.97"
$ws.Cells.Item(2, 5).Value = "  +0.73%  "

$ws.Cells.Item(3, 4).Value = "1.641.55"
$ws.Cells.Item(3, 5).Value = "  +1.23%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "215.81"
$ws.Cells.Item(5, 5).Value = "  +0.68%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.5082"
$ws.Cells.Item(6, 5).Value = "  +0.25%  "

$ws.Cells.Item(7, 5).Value = "  +0.13%  "

$ws.Cells.Item(8, 5).Value = "  +1.81%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06466"
$ws.Cells.Item(9, 5).Value = "  +1.41%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "20.25"
$ws.Cells.Item(10, 5).Value = "  +5.31%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07811"
$ws.Cells.Item(11, 5).Value = "  +0.51%  "

$ws.Cells.Item(12, 4).Value = "1.658.46"
$ws.Cells.Item(12, 5).Value = "  +2.29%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.269"
$ws.Cells.Item(13, 5).Value = "  +0.89%  "

$ws.Cells.Item(14, 4).Value = "1.867.35"
$ws.Cells.Item(14, 5).Value = "  +1.20%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.5663"

$ws.Cells.Item(16, 4).Value = "0.0₅7716"
$ws.Cells.Item(16, 5).Value = "  +2.63%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "63.57"
$ws.Cells.Item(17, 5).Value = "  +0.19%  "

$ws.Cells.Item(18, 4).Value = "25.920.32"
$ws.Cells.Item(18, 5).Value = "  +0.69%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.003"
$ws.Cells.Item(19, 5).Value = "  +0.13%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "194.92"
$ws.Cells.Item(20, 5).Value = "  +0.91%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "4.403"
$ws.Cells.Item(21, 5).Value = "  +0.66%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "9.984"
$ws.Cells.Item(22, 5).Value = "  +2.53%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.241"
$ws.Cells.Item(23, 5).Value = "  +4.99%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.003"
$ws.Cells.Item(24, 5).Value = "  +0.09%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.763"
$ws.Cells.Item(25, 5).Value = "  -5.30%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "138.63"
$ws.Cells.Item(26, 5).Value = "  -1.30%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.1231"
$ws.Cells.Item(27, 5).Value = "  -1.33%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.872"
$ws.Cells.Item(28, 5).Value = "  +2.37%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.60"
$ws.Cells.Item(29, 5).Value = "  +1.05%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.245"
$ws.Cells.Item(30, 5).Value = "  +1.02%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.05027"
$ws.Cells.Item(31, 5).Value = "  +3.57%  "

$ws.Cells.Item(32, 5).Value = "  +0.63%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.259"
$ws.Cells.Item(33, 5).Value = "  +2.90%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.581"
$ws.Cells.Item(34, 5).Value = "  +2.78%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.384"
$ws.Cells.Item(35, 5).Value = "  +0.97%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.9104"
$ws.Cells.Item(36, 5).Value = "  +2.30%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.586"
$ws.Cells.Item(37, 5).Value = "  +2.14%  "

$ws.Cells.Item(38, 5).Value = "  +1.14%  "

$ws.Cells.Item(39, 4).Value = "1.128.86"
$ws.Cells.Item(39, 5).Value = "  +0.49%  "

$ws.Cells.Item(40, 5).Value = "  +1.31%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.9940"
$ws.Cells.Item(41, 5).Value = "  -0.66%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "100.03"
$ws.Cells.Item(42, 5).Value = "  +3.20%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.504"
$ws.Cells.Item(43, 5).Value = "  -1.02%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.8021"
$ws.Cells.Item(44, 5).Value = "  +1.15%  "

$ws.Cells.Item(45, 4).Value = "0.0₈110"
$ws.Cells.Item(45, 5).Value = "  -3.60%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "55.77"
$ws.Cells.Item(46, 5).Value = "  +2.29%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.4235"
$ws.Cells.Item(47, 5).Value = "  -4.07%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.698"
$ws.Cells.Item(48, 5).Value = "  +1.88%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.05045"
$ws.Cells.Item(49, 5).Value = "  -1.43%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.003"
$ws.Cells.Item(50, 5).Value = "  +0.64%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.002"
$ws.Cells.Item(51, 5).Value = "  +0.21%  "
